$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: write a cell as TEXT using the "new" plain text style (no wrap,
# text number format) that matches cellXfs index 4 in the target workbook
# (numFmtId 49 "@" with no alignment/wrap). Resetting the style to "Normal"
# before re-applying the text number format strips any inherited column
# formatting (e.g. wrapText) so the freshly created style stays clean.
# ---------------------------------------------------------------------------
function Set-TextCell($addr, $value) {
    $r = $ws.Range($addr)
    $r.Style = "Normal"
    $r.NumberFormat = "@"
    $r.Value = $value
}

# Cells that already carry the (untouched) "hyperlink-like" style index 2 -
# only their values change, the style is left exactly as-is.
function Set-ValueOnly($addr, $value) {
    $ws.Range($addr).Value = $value
}

# Cells that must end up with NO explicit style at all (default style 0) -
# reset to Normal (General number format) and write the value.
function Set-PlainCell($addr, $value) {
    $r = $ws.Range($addr)
    $r.Style = "Normal"
    $r.Value = $value
}

# ----------------------------- Row 2 ---------------------------------------
Set-TextCell "B2" "29.00"
Set-TextCell "C2" "29"
Set-TextCell "D2" "Temperatues are within Variance Range"
Set-TextCell "E2" "Mist and mist"
Set-TextCell "F2" "Humid and Overcast`n"
Set-TextCell "G2" "Both portals show slightly similar Weather conditions!"
Set-TextCell "H2" "83"
Set-TextCell "I2" "74"
Set-ValueOnly "K2" "1.50"
Set-TextCell "L2" "3.3949999999999996"

# ----------------------------- Row 3 ---------------------------------------
Set-TextCell "B3" "28.57"
Set-ValueOnly "C3" "27"
Set-ValueOnly "D3" "Temperatues are within Variance Range"
Set-TextCell "E3" "Rain and light rain"
Set-TextCell "F3" "Humid and Partly Cloudy`n"
Set-TextCell "G3" "Both portals show different Weather conditions!"
Set-TextCell "H3" "73"
Set-TextCell "I3" "83"
Set-ValueOnly "K3" "5.10"
Set-TextCell "L3" "5.75"

# ----------------------------- Row 4 ---------------------------------------
Set-TextCell "B4" "30.07"
Set-ValueOnly "C4" "29"
Set-ValueOnly "D4" "Temperatues are within Variance Range"
Set-TextCell "E4" "Clouds and broken clouds"
Set-TextCell "F4" "Humid and Mostly Cloudy`n"
Set-TextCell "G4" "Both portals show slightly similar Weather conditions!"
Set-TextCell "H4" "66"
Set-TextCell "I4" "77"
Set-ValueOnly "K4" "7.04"
Set-TextCell "L4" "7.49"

# ----------------------------- Row 5 ---------------------------------------
Set-TextCell "B5" "24.00"
Set-TextCell "C5" "24"
Set-TextCell "D5" "Temperatues are within Variance Range"
Set-TextCell "E5" "Mist and mist"
Set-TextCell "F5" "Overcast`n"
Set-PlainCell "G5" "Both portals show different Weather conditions!"
Set-TextCell "H5" "88"
Set-TextCell "I5" "85"
Set-TextCell "K5" "1.50"
Set-TextCell "L5" "3.2199999999999998"

# ----------------------------- Row 6 ---------------------------------------
Set-TextCell "B6" "29.00"
Set-TextCell "C6" "28"
Set-TextCell "D6" "Temperatues are within Variance Range"
Set-TextCell "E6" "Haze and haze"
Set-TextCell "F6" "Humid and Mostly Cloudy`n"
Set-PlainCell "G6" "Both portals show slightly similar Weather conditions!"
Set-TextCell "H6" "89"
Set-TextCell "I6" "94"
Set-TextCell "K6" "1.50"
Set-TextCell "L6" "5.53"

# ----------------------------- Row 7 ---------------------------------------
Set-TextCell "B7" "32.00"
Set-TextCell "C7" "31"
Set-TextCell "D7" "Temperatues are within Variance Range"
Set-TextCell "E7" "Haze and haze"
Set-TextCell "F7" "Humid and Overcast`n"
Set-PlainCell "G7" "Both portals show slightly similar Weather conditions!"
Set-TextCell "H7" "84"
Set-TextCell "I7" "72"
Set-TextCell "K7" "5.10"
Set-TextCell "L7" "7.41"

# ------------------------- Column width tweaks ------------------------------
# Columns D and G had their "best fit" widths recomputed once real data (and,
# for G, new content) was in place. Nudge them to the closest width this
# engine's column-width model can express.
$ws.Columns("D:D").ColumnWidth = 36
$ws.Columns("G:G").ColumnWidth = 48.5

# ------------------------------ Selection -----------------------------------
$ws.Range("G9").Select()
